$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "Ineligible" column values (S, T, U, V) for rows 3-7 ---
$ws.Range("S3").Value = 222
$ws.Range("T3").Value = 600

$ws.Range("S4").Value = 100
$ws.Range("U4").Value = 744

$ws.Range("S5").Value = 333
$ws.Range("T5").Value = 450

$ws.Range("T6").Value = 600
$ws.Range("V6").Value = 700

$ws.Range("S7").Value = 560
$ws.Range("U7").Value = 766
$ws.Range("V7").Value = 987

# --- Re-touch the header merged cells (unmerge + remerge) so the stored
#     mergeCells order matches how Excel rewrote the sheet ---
$mergeOrder = @("K1:K2","O1:R1","S1:V1","N1:N2","M1:M2","L1:L2","J1:J2","I1:I2","H1:H2","G1:G2","A1:A2","B1:B2","C1:C2","D1:D2","E1:E2","F1:F2","W1:W2","X1:X2","Y1:Y2","Z1:Z2","AA1:AA2")
foreach ($r in $mergeOrder) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $mergeOrder) {
    $ws.Range($r).Merge()
}

# --- Update the current selection on the sheet ---
$ws.Range("V8").Select()
